# Update gh-pages generated output figures (想去人数 / 最低票价)
# for the "合肥·第九届环形宇宙动漫游戏嘉年华" event (id=92565),
# and the "想去人数" for "合肥·W·A第五人格同人only2.0" (id=91123),
# on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 414
$ws1.Range("F3").Value = 2459
$ws1.Range("G3").Value = 72

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 414
$ws4.Range("F7").Value = 2459
$ws4.Range("G7").Value = 72
